$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '64.062.35'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '3.153.36'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.02'
$ws.Range("E5").Value = '  -1.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.81'
$ws.Range("E6").Value = '  -2.82%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.148.29'
$ws.Range("E8").Value = '  -0.58%  '
$ws.Range("E9").Value = '  -0.99%  '
$ws.Range("E10").Value = '  -1.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.38'
$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.467'
$ws.Range("E12").Value = '  -1.58%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.02'
$ws.Range("E14").Value = '  -2.39%  '
$ws.Range("D15").Value = '3.675.37'
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("D17").Value = '64.136.55'
$ws.Range("E17").Value = '  -0.83%  '
$ws.Range("D18").Value = '3.155.16'
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("E19").Value = '  -1.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '493.34'
$ws.Range("E20").Value = '  +2.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.75'
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.711'
$ws.Range("E22").Value = '  -1.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.66'
$ws.Range("E23").Value = '  -4.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.83'
$ws.Range("E24").Value = '  +3.92%  '
$ws.Range("E25").Value = '  -3.30%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  -2.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.22'
$ws.Range("E28").Value = '  -4.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.00'
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.06'
$ws.Range("E30").Value = '  -1.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.57'
$ws.Range("E31").Value = '  +3.69%  '
$ws.Range("E32").Value = '  -5.37%  '
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("E34").Value = '  -2.04%  '
$ws.Range("E35").Value = '  -2.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.03'
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.74'
$ws.Range("E37").Value = '  -0.53%  '
$ws.Range("D38").Value = '0.0₃0744'
$ws.Range("E38").Value = '  -5.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.97'
$ws.Range("E39").Value = '  -7.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '435.71'
$ws.Range("E40").Value = '  -5.60%  '
$ws.Range("E41").Value = '  -0.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.119'
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.31'
$ws.Range("E43").Value = '  -1.49%  '
$ws.Range("D44").Value = '2.944.07'
$ws.Range("E44").Value = '  +3.24%  '
$ws.Range("E45").Value = '  -3.91%  '
$ws.Range("E46").Value = '  -5.94%  '
$ws.Range("E47").Value = '  -2.59%  '
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.97'
$ws.Range("E49").Value = '  -2.71%  '
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("E51").Value = '  +0.42%  '
